$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure column D cells are treated as Text so numeric-looking strings
# (e.g. "0.999", "6.68") are not auto-converted to Number values by the
# Value setter -- matches the source file where every cell is inline text.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = '59.381.88'
$ws.Range("E2").Value = '  +0.74%  '

$ws.Range("D3").Value = '2.589.03'
$ws.Range("E3").Value = '  -0.48%  '

$ws.Range("D5").Value = '570.81'
$ws.Range("E5").Value = '  +3.08%  '

$ws.Range("D6").Value = '143.96'
$ws.Range("E6").Value = '  +0.05%  '

$ws.Range("D7").Value = '0.999'
$ws.Range("E7").Value = '  +0.11%  '

$ws.Range("D8").Value = '0.597'
$ws.Range("E8").Value = '  +0.23%  '

$ws.Range("D9").Value = '2.599.44'
$ws.Range("E9").Value = '  -0.40%  '

$ws.Range("D10").Value = '6.68'
$ws.Range("E10").Value = '  -1.55%  '

$ws.Range("D11").Value = '0.104'
$ws.Range("E11").Value = '  +3.60%  '

$ws.Range("D12").Value = '0.159'
$ws.Range("E12").Value = '  +11.32%  '

$ws.Range("D13").Value = '0.345'
$ws.Range("E13").Value = '  +2.84%  '

$ws.Range("D14").Value = '3.049.01'
$ws.Range("E14").Value = '  -0.30%  '

$ws.Range("D15").Value = '59.362.13'
$ws.Range("E15").Value = '  +0.78%  '

$ws.Range("D16").Value = '22.58'
$ws.Range("E16").Value = '  +7.85%  '

$ws.Range("D17").Value = '0.0000137'
$ws.Range("E17").Value = '  +3.94%  '

$ws.Range("D18").Value = '2.595.13'
$ws.Range("E18").Value = '  -0.49%  '

$ws.Range("D19").Value = '4.54'
$ws.Range("E19").Value = '  +1.32%  '

$ws.Range("D20").Value = '335.67'
$ws.Range("E20").Value = '  -0.76%  '

$ws.Range("D21").Value = '10.28'
$ws.Range("E21").Value = '  +1.51%  '

$ws.Range("D22").Value = '6.21'
$ws.Range("E22").Value = '  +0.61%  '

$ws.Range("E23").Value = '  -0.08%  '

$ws.Range("D24").Value = '64.51'
$ws.Range("E24").Value = '  -3.33%  '

$ws.Range("D25").Value = '0.458'
$ws.Range("E25").Value = '  +6.80%  '

$ws.Range("B26").Value = 'Binance-PegBSC-USD'
$ws.Range("C26").Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range("D26").Value = '0.992'
$ws.Range("E26").Value = '  -0.44%  '

$ws.Range("B27").Value = 'Kaspa'
$ws.Range("C27").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D27").Value = '0.161'
$ws.Range("E27").Value = '  -0.22%  '

$ws.Range("D28").Value = '7.30'
$ws.Range("E28").Value = '  +1.49%  '

$ws.Range("D29").Value = '0.0₃0783'
$ws.Range("E29").Value = '  +3.47%  '

$ws.Range("E30").Value = '  +0.02%  '

$ws.Range("D31").Value = '1.68'
$ws.Range("E31").Value = '  +0.41%  '

$ws.Range("D32").Value = '6.11'
$ws.Range("E32").Value = '  +1.32%  '

$ws.Range("D33").Value = '158.64'
$ws.Range("E33").Value = '  +2.92%  '

$ws.Range("D34").Value = '19.07'
$ws.Range("E34").Value = '  +0.38%  '

$ws.Range("D35").Value = '4.07'
$ws.Range("E35").Value = '  +3.27%  '

$ws.Range("D36").Value = '1.15'
$ws.Range("E36").Value = '  +1.67%  '

$ws.Range("B37").Value = 'SuiNetwork'
$ws.Range("C37").Value = 'https://coinranking.com/coin/3xJluUMvp+suinetwork-sui'
$ws.Range("D37").Value = '0.877'
$ws.Range("E37").Value = '  -2.47%  '

$ws.Range("B38").Value = 'Fetch.AI'
$ws.Range("C38").Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range("D38").Value = '0.882'
$ws.Range("E38").Value = '  +0.13%  '

$ws.Range("D39").Value = '37.16'
$ws.Range("E39").Value = '  +0.59%  '

$ws.Range("E40").Value = '  +2.05%  '

$ws.Range("D41").Value = '294.89'
$ws.Range("E41").Value = '  +4.12%  '

$ws.Range("D42").Value = '3.69'
$ws.Range("E42").Value = '  +2.12%  '

$ws.Range("E43").Value = '  +0.13%  '

$ws.Range("E44").Value = '  +2.63%  '

$ws.Range("E45").Value = '  -0.80%  '

$ws.Range("B46").Value = 'Hedera'
$ws.Range("C46").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D46").Value = '0.0540'
$ws.Range("E46").Value = '  +0.59%  '

$ws.Range("B47").Value = 'EnergySwap'
$ws.Range("C47").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D47").Value = '19.34'
$ws.Range("E47").Value = '  +2.60%  '

$ws.Range("D48").Value = '10.64'
$ws.Range("E48").Value = '  +0.00%  '

$ws.Range("D49").Value = '125.77'
$ws.Range("E49").Value = '  +6.88%  '

$ws.Range("D50").Value = '0.0233'
$ws.Range("E50").Value = '  +2.05%  '

$ws.Range("D51").Value = '18.62'
$ws.Range("E51").Value = '  +2.76%  '

# Restore the default (unstyled) cell style now that the text values are
# set, so column D does not pick up a persistent "@" number format that
# was never part of the original formatting.
$ws.Range("D2:D51").Style = "Normal"
